$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 161, shifting existing rows 161:259 down to 162:260
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new record
$ws.Range("A161").Value = 7
$ws.Range("B161").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C161").Value = "Ñuble"
$ws.Range("D161").Value = 44830
$ws.Range("E161").Value = 16
$ws.Range("F161").Value = 100112043
$ws.Range("G161").Value = "Pepino ensalada"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 100
$ws.Range("K161").Value = 19000
$ws.Range("L161").Value = 20000
$ws.Range("M161").Value = 19500
$ws.Range("N161").Value = "$/caja 60 unidades"
$ws.Range("O161").Value = "Región de Arica y Parinacota"
$ws.Range("P161").Value = 325
$ws.Range("Q161").Value = 60
$ws.Range("R161").Value = "Hortaliza"
